$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data for "Just greater number" (row 7)
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = "Binary S"
$ws.Range("D7").Value = "Just greater number"
$ws.Range("E7").Value = "No Link"

$dot = [char]0x00B7
$f7Text = "dsa/5_just_greater_number.java at main " + $dot + " ankurnecessary/dsa " + $dot + " GitHub"
$ws.Range("F7").Value = $f7Text

# Add the hyperlink for F7 (this also applies a Hyperlink-like cell style)
$ws.Hyperlinks.Add(
    $ws.Range("F7"),
    "https://github.com/ankurnecessary/dsa/blob/main/1_binarySearch/5_just_greater_number.java",
    "",
    "",
    $f7Text
) | Out-Null

# Re-apply the exact formatting used by the other Github-link cells (F4:F6)
# so F7 shares the same cell style rather than Excel minting a near-duplicate one.
$ws.Range("F6").Copy() | Out-Null
$ws.Range("F7").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Match target row height (auto height for the shorter "No Link" row)
$ws.Rows.Item(7).RowHeight = 43.2

# Select F7 as the final active cell, matching the saved selection state
$ws.Range("F7").Select() | Out-Null
